$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeiterfassung")

# Duplicate formatting of the last existing booking row (13) into the new row (14)
$ws.Range("A13:D13").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New booking entry (row 14)
$ws.Range("A14").Value = 45586
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Planung"
$ws.Range("D14").Value = "Mail an Management & Teams Aufgabenplaner aktualisiert."

$ws.Range("A15").Select()
